$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.053.87"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "3.099.99"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D5").Value = "523.16"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "144.02"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "7.39"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "3.635.41"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "27.00"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "59.027.61"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("D17").Value = "3.105.78"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "6.17"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "343.88"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "0.508"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").Value = "65.79"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "0.0₃0930"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").Value = "6.74"
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("D32").Value = "21.00"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "155.06"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "4.64"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "6.16"
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("D36").Value = "26.78"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("D38").Value = "0.0687"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "3.96"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "3.145.62"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").Value = "36.78"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.665"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("D45").Value = "2.300.08"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").Value = "0.0256"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "20.90"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").Value = "0.963"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "6.02"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "0.755"
$ws.Range("E50").Value = "  +9.49%  "
$ws.Range("D51").Value = "264.83"
$ws.Range("E51").Value = "  +12.61%  "
